# Apply the "calorimetry : scripts : tests : updated" edit:
#  1. Rename the "adj_r_squared" sheet to "metrics".
#  2. Replace its single R^2 cell with a small name/value metrics table
#     (Adj.R^2, NRMSE, SMAPE, RMSE).
#  3. Drop the two leftover blank/zero rows from "input_enthalpies"
#     (sheet keeps only its header row).

$wb = $excel.ActiveWorkbook

# --- 1 & 2: adj_r_squared -> metrics, new metrics table -------------------
$metrics = $wb.Worksheets.Item("adj_r_squared")
$metrics.Name = "metrics"

$metrics.Range("A1").Value = "metrics"
$metrics.Range("B1").Value = "value"

$metrics.Range("A2").Value = "Adj.R^2"
$metrics.Range("B2").Value = 0.996023014013136

$metrics.Range("A3").Value = "NRMSE"
$metrics.Range("B3").Value = 0.0564619206417458

$metrics.Range("A4").Value = "SMAPE"
$metrics.Range("B4").Value = 0.319540708155668

$metrics.Range("A5").Value = "RMSE"
$metrics.Range("B5").Value = 0.00215614789773029

# --- 3: input_enthalpies - remove the extra rows ---------------------------
$enthalpies = $wb.Worksheets.Item("input_enthalpies")
$enthalpies.Range("A2:B3").Delete()
